$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the existing table ("right_arm") from A1:F5 to A1:H5 so two new
# columns (Dir / Dir ) can be inserted: one before "Mid" and one after "Max".
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:H5"))

# ---- Header row -----------------------------------------------------
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Num"
$ws.Range("C1").Value = "Degree (A)"
$ws.Range("D1").Value = "Min"
$ws.Range("E1").Value = "Dir"
$ws.Range("H1").Value = "Dir "
$ws.Range("F1").Value = "Mid"
$ws.Range("G1").Value = "Max"

# ---- Row 2: right_shoulder_x / M4 ------------------------------------
$ws.Range("A2").Value = "right_shoulder_x"
$ws.Range("B2").Value = "M4"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("H2").Value = "Inside"
$ws.Range("E2").Value = "Outside"
$ws.Range("F2").Value = 130
$ws.Range("G2").Value = 180

# ---- Row 3: right_shoulder_y / M3 ------------------------------------
$ws.Range("A3").Value = "right_shoulder_y"
$ws.Range("B3").Value = "M3"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = "Backward"
$ws.Range("H3").Value = "Forward"
$ws.Range("F3").Value = 150
$ws.Range("G3").Value = 180

# ---- Row 4: right_shoulder_z / M1 ------------------------------------
$ws.Range("A4").Value = "right_shoulder_z"
$ws.Range("B4").Value = "M1"
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 127
$ws.Range("E4").Value = "Inside"
$ws.Range("H4").Value = "Outside"
$ws.Range("F4").Value = 130
$ws.Range("G4").Value = 135

# ---- Row 5: right_bicept / M2 -----------------------------------------
$ws.Range("A5").Value = "right_bicept"
$ws.Range("B5").Value = "M2"
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 70
$ws.Range("E5").Value = "down"
$ws.Range("H5").Value = "up"
$ws.Range("F5").Value = 80
$ws.Range("G5").Value = 173

# Move the active cell / selection the way the saved workbook shows it.
[void]$ws.Range("F10").Select()

# Widen the workbook window (best-effort; some hosts ignore this).
try { $excel.ActiveWindow.Width = 21600 } catch {}
try { $wb.Windows.Item(1).Width = 21600 } catch {}
